# Remove the extra "Context" / "element:Medication" row from the Metadata
# sheet (row 22 duplicated the Context property; only the
# element:Medication.ingredient row should remain) and clear the stray
# Pattern value ("true") on the Extension.value[x] row of the Elements
# sheet.

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$meta.Rows.Item(21).Delete()
$meta.Range("B8").Value = "2025-04-09T16:24:06+00:00"

$elements = $wb.Worksheets.Item("Elements")
$elements.Range("S6").Value = "'"
$elements.Range("D2").Copy()
$elements.Range("S6").PasteSpecial(-4122)
